# Update PM04 Tidsregistrering for Bille.xlsx
# Adds four new time-registration entries (rows 11-14) on "Ark1",
# and updates the view/selection state on both sheets.

$wb = $excel.ActiveWorkbook

# ---- Ark1: fill in the four new time-registration rows ------------------
$ws1 = $wb.Worksheets.Item("Ark1")

# Row 11 - OC0802 / System Analyst
$ws1.Range("A11").Value = "OC0802"
$ws1.Range("B11").Value = "System Analyst "
$ws1.Range("C4").Copy()
$ws1.Range("C11").PasteSpecial(-4122)
$ws1.Range("C11").Value = 43887
$ws1.Range("D11").Value = 0.39583333333333331
$ws1.Range("E11").Value = 0.44097222222222227

# Row 12 - KKO use-case moede / business-Process Analyst
$ws1.Range("A12").Value = "KKO use-case møde"
$ws1.Range("B12").Value = "business-Process Analyst"
$ws1.Range("C4").Copy()
$ws1.Range("C12").PasteSpecial(-4122)
$ws1.Range("C12").Value = 43887
$ws1.Range("D12").Value = 0.46875
$ws1.Range("E12").Value = 0.49652777777777773

# Row 13 - IndtjeningsBidrag AD moede / business-Process Analyst
$ws1.Range("A13").Value = "IndtjeningsBidrag AD møde"
$ws1.Range("B13").Value = "business-Process Analyst"
$ws1.Range("C4").Copy()
$ws1.Range("C13").PasteSpecial(-4122)
$ws1.Range("C13").Value = 43887
$ws1.Range("D13").Value = 0.51041666666666663
$ws1.Range("E13").Value = 0.54513888888888895

# Row 14 - DOM06 KontantKapacitetsOmkostning / business-Process Analyst
$ws1.Range("A14").Value = "DOM06 KontantKapacitetsOmkostning"
$ws1.Range("B14").Value = "business-Process Analyst"
$ws1.Range("C4").Copy()
$ws1.Range("C14").PasteSpecial(-4122)
$ws1.Range("C14").Value = 43887
$ws1.Range("D14").Value = 0.54166666666666663
$ws1.Range("E14").Value = 0.59375

$excel.CutCopyMode = 0

# ---- Ark2: refresh its view state (drops stale scroll position) --------
$ws2 = $wb.Worksheets.Item("Ark2")
$ws2.Activate()
$ws2.Range("B5").Select()

# ---- Ark1: make it the active sheet again and park the selection -------
$ws1.Activate()
$ws1.Range("B13").Select()
